$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 25.02.2022 14:00"

# Row 4 (Globus): new price found is 37.9, old price (previous B4) becomes C4
$ws.Range("C4").Value = $ws.Range("B4").Value2
$ws.Range("B4").Value = 37.9

# Delta column becomes a formatted text string rather than a numeric delta
$ws.Range("D4").Value = "'+1.0"

# Old Datum column becomes a literal text timestamp (no date number format)
$ws.Range("E4").Style = "Normal"
$ws.Range("E4").Value = "2022-02-25 14:00:11"
